# Database model.docx edit:
#  - "High score table" section becomes "Games" with new sub-items,
#    followed by new "Levels" and "Maps" sections.
#
# Strategy: capture two structural paragraph "templates" (one List
# Paragraph-styled, one Normal-styled) from the existing document while
# their formatting is still pristine, then use Range.InsertAfter() with
# the template's text (which already ends in the paragraph-mark \r) to
# stamp out new paragraphs with the right pPr/numPr. Range.InsertAfter
# is anchored on the END of a paragraph that is NOT the very last
# paragraph in the document body (Word's COM model special-cases a
# zero-length range sitting at the final paragraph mark), so a throw-away
# paragraph is appended first to push the real content off that boundary.

$d = $word.ActiveDocument

# Templates captured before we mutate anything.
$listTemplate   = $d.Paragraphs(11).Range.Text   # "High score ID..." - List Paragraph, numId 1
$normalTemplate = $d.Paragraphs(10).Range.Text   # "High score table" - Normal

# --- "High score table" heading -> "Games" -------------------------------
$d.Paragraphs(10).Range.Text = "Games"

# --- "High score ID..." -> "Starts" ---------------------------------------
$d.Paragraphs(11).Range.Text = "Starts"

# --- "Score (= user['high score'])" (multi-run) -> "Ends" -----------------
$p12 = $d.Paragraphs(12)
$p12s = $p12.Range.Start
$p12e = $p12.Range.End - 1
$d.Range($p12s, $p12e).Text = "Ends"

# --- insert "Game type id" / "Player id" list items after "Ends" ----------
$anchor = $d.Paragraphs(12).Range.End
$r = $d.Range($anchor, $anchor)
$r.InsertAfter($listTemplate)
$r.InsertAfter($listTemplate)
$d.Paragraphs(13).Range.Text = "Game type id"
$d.Paragraphs(14).Range.Text = "Player id"

# --- "User ID" (+ bookmark) -> "Result id" ---------------------------------
$p15 = $d.Paragraphs(15)
$p15s = $p15.Range.Start
$p15e = $p15.Range.End - 1
$d.Range($p15s, $p15e).Text = "Result id"

# --- push a throw-away paragraph to the very end so the next InsertAfter
#     anchor (end of the "Result id" paragraph) is no longer the last
#     paragraph mark in the document -----------------------------------
$d.Paragraphs(15).Range.InsertParagraphAfter()

# --- insert the "Levels" and "Maps" sections after "Result id" ------------
$anchor2 = $d.Paragraphs(15).Range.End
$r2 = $d.Range($anchor2, $anchor2)

$r2.InsertAfter($normalTemplate)   # -> 16: empty (Normal)
$r2.InsertAfter($normalTemplate)   # -> 17: "Levels "
$r2.InsertAfter($listTemplate)     # -> 18: "Level id"
$r2.InsertAfter($listTemplate)     # -> 19: "Time limit"
$r2.InsertAfter($normalTemplate)   # -> 20: empty (Normal)
$r2.InsertAfter($normalTemplate)   # -> 21: "Maps"
$r2.InsertAfter($listTemplate)     # -> 22: "Map id"
$r2.InsertAfter($listTemplate)     # -> 23: "Level id"
# 24 is the throw-away paragraph pushed to the end earlier; reuse it for
# the final "High score" list item instead of leaving an extra blank one.

$d.Paragraphs(16).Range.Text = ""
$d.Paragraphs(17).Range.Text = "Levels "
$d.Paragraphs(18).Range.Text = "Level id"
$d.Paragraphs(19).Range.Text = "Time limit"
$d.Paragraphs(20).Range.Text = ""
$d.Paragraphs(21).Range.Text = "Maps"
$d.Paragraphs(22).Range.Text = "Map id"
$d.Paragraphs(23).Range.Text = "Level id"

$last = $d.Paragraphs(24)
$last.Style = "List Paragraph"
$last.Range.ListFormat.ApplyListTemplateWithLevel($d.Paragraphs(23).Range.ListFormat.ListTemplate)
$last.Range.Text = "High score"

Write-Host "Final paragraph count:" $d.Paragraphs.Count
foreach ($pp in $d.Paragraphs) {
    Write-Host "[" $pp.Range.Text "]"
}
